$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (bold font, border, alignment) of the existing "IP" header cell
# onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data for columns I (I0) and J (IF), rows 2-20
$data = @{
    2  = @(8, 9)
    3  = @(9, 9)
    4  = @(9, 9)
    5  = @(8, 9)
    6  = @(6, 6)
    7  = @(9, 9)
    8  = @(6, 7)
    9  = @(8, 8)
    10 = @(4, 5)
    11 = @(7, 7)
    12 = @(6, 6)
    13 = @(6, 6)
    14 = @(8, 8)
    15 = @(6, 6)
    16 = @(6, 6)
    17 = @(7, 7)
    18 = @(7, 8)
    19 = @(5, 5)
    20 = @(9, 9)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
